# Apply the commit's changes:
#  - Rename the two "biotop" category labels used in the samples sheet:
#      "coniferous forest (regular)"            -> "regular forest"
#      "coniferous forest (post-fire succession)" -> "post-fire forest"
#  - Switch the active tab from "main" to "samples"
#  - Update the selection on the "samples" sheet to C9

$wb = $excel.ActiveWorkbook

$samplesWs = $wb.Worksheets.Item("samples")

# Rename the biotop strings on the samples sheet (column C, rows 2-57).
# The new shared-string table entries get created in first-use order, so
# touch the "post-fire" rows before the "regular" rows to mirror the
# upstream edit (post-fire forest ends up before regular forest in sst).
for ($r = 2; $r -le 57; $r++) {
    $cell = $samplesWs.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -eq "coniferous forest (post-fire succession)") {
        $cell.Value = "post-fire forest"
    }
}
for ($r = 2; $r -le 57; $r++) {
    $cell = $samplesWs.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -eq "coniferous forest (regular)") {
        $cell.Value = "regular forest"
    }
}

# Switch the active sheet to "samples" and move the selection to C9
$samplesWs.Activate()
$samplesWs.Range("C9").Select()
